# Update the cryptos price/volume table (columns D = Price, E = Volume(1h)).
# Price values that look like plain decimals (e.g. "394.69") are prefixed
# with a leading apostrophe so Excel stores them as literal text instead of
# coercing them to numbers (matching the original inline-string cells, and
# preserving formatting like trailing zeros, e.g. "74.00"/"2.10").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.026.08'
$ws.Range("E2").Value = '  +6.71%  '
$ws.Range("D3").Value = '3.239.32'
$ws.Range("E3").Value = '  +3.23%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '''394.69'
$ws.Range("E5").Value = '  -0.48%  '
$ws.Range("D6").Value = '''107.11'
$ws.Range("E6").Value = '  -2.11%  '
$ws.Range("D7").Value = '''0.574'
$ws.Range("E7").Value = '  +5.12%  '
$ws.Range("D8").Value = '3.235.54'
$ws.Range("E8").Value = '  +3.33%  '
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").Value = '''0.618'
$ws.Range("E10").Value = '  +1.46%  '
$ws.Range("D11").Value = '''39.03'
$ws.Range("E11").Value = '  +0.47%  '
$ws.Range("D12").Value = '''0.0959'
$ws.Range("E12").Value = '  +10.13%  '
$ws.Range("E13").Value = '  +1.65%  '
$ws.Range("D14").Value = '3.748.11'
$ws.Range("E14").Value = '  +2.66%  '
$ws.Range("D15").Value = '''8.18'
$ws.Range("E15").Value = '  +2.31%  '
$ws.Range("D16").Value = '''19.13'
$ws.Range("E16").Value = '  +0.36%  '
$ws.Range("D17").Value = '3.233.74'
$ws.Range("E17").Value = '  +2.67%  '
$ws.Range("E18").Value = '  -2.04%  '
$ws.Range("D19").Value = '''10.88'
$ws.Range("E19").Value = '  +3.49%  '
$ws.Range("D20").Value = '56.852.46'
$ws.Range("E20").Value = '  +6.33%  '
$ws.Range("E21").Value = '  +2.08%  '
$ws.Range("E22").Value = '  +8.48%  '
$ws.Range("E23").Value = '  +1.58%  '
$ws.Range("D24").Value = '''296.38'
$ws.Range("E24").Value = '  +9.42%  '
$ws.Range("D25").Value = '''74.00'
$ws.Range("E25").Value = '  +4.30%  '
$ws.Range("E26").Value = '  -2.24%  '
$ws.Range("D27").Value = '''27.82'
$ws.Range("E27").Value = '  +1.34%  '
$ws.Range("E28").Value = '  -3.72%  '
$ws.Range("D29").Value = '''7.27'
$ws.Range("E29").Value = '  -0.42%  '
$ws.Range("D30").Value = '''0.168'
$ws.Range("E30").Value = '  -0.67%  '
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("D32").Value = '''11.40'
$ws.Range("E32").Value = '  +4.00%  '
$ws.Range("E33").Value = '  -1.14%  '
$ws.Range("D34").Value = '''37.91'
$ws.Range("E34").Value = '  +1.95%  '
$ws.Range("D35").Value = '''0.0483'
$ws.Range("E35").Value = '  -3.29%  '
$ws.Range("D36").Value = '''2.11'
$ws.Range("E36").Value = '  +1.60%  '
$ws.Range("D37").Value = '''51.78'
$ws.Range("E37").Value = '  +2.76%  '
$ws.Range("E38").Value = '  -3.21%  '
$ws.Range("E39").Value = '  -0.17%  '
$ws.Range("D40").Value = '''2.95'
$ws.Range("E40").Value = '  +5.78%  '
$ws.Range("D41").Value = '''134.36'
$ws.Range("E41").Value = '  +3.38%  '
$ws.Range("E42").Value = '  +2.64%  '
$ws.Range("D43").Value = '''3.96'
$ws.Range("E43").Value = '  -3.26%  '
$ws.Range("E44").Value = '  -0.14%  '
$ws.Range("D45").Value = '''17.02'
$ws.Range("E45").Value = '  -1.16%  '
$ws.Range("E46").Value = '  -2.62%  '
$ws.Range("D47").Value = '''22.13'
$ws.Range("E47").Value = '  +0.06%  '
$ws.Range("D48").Value = '2.157.94'
$ws.Range("E48").Value = '  +3.87%  '
$ws.Range("D49").Value = '''2.10'
$ws.Range("E49").Value = '  +1.43%  '
$ws.Range("D50").Value = '''2.01'
$ws.Range("E50").Value = '  +21.61%  '
$ws.Range("E51").Value = '  -2.77%  '
